$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 9 (shifts existing rows 9..40 down to 10..41)
$ws.Rows.Item(9).Insert()

# Populate the new row 9 with the "Alcohol en gel" product
$ws.Range("A9").Value = 7791274196522
$ws.Range("B9").Value = "Alcohol"
$ws.Range("C9").Value = "para manos"
$ws.Range("D9").Value = "en gel"
$ws.Range("E9").Value = "Algabo"
$ws.Range("F9").Value = 300
$ws.Range("G9").Value = "ml."
$ws.Range("H9").Value = "Botella"
$ws.Range("I9").Value = "Alcoholes"
$ws.Range("J9").Value = "Argentina"
$ws.Range("K9").Value = 6
$ws.Range("L9").Value = $false
$ws.Range("M9").Value = $true
$ws.Range("O9").Value = $true

# O9 inherited the plain "copied from row above" look on Insert; restore the
# thin gray boxed border used by the rest of the ImagenExactaDelArticulo column
foreach ($edge in 7,8,9,10) {
    $ws.Range("O9").Borders.Item($edge).LineStyle = 1
    $ws.Range("O9").Borders.Item($edge).Color = 12632256
}


